$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "Dimensional statements. Occurrences contexts: Events. Dimensional
# contexts: Context occurrences hierarchy for order relations
# assertions." ->
# "Dimensional contexts: Contexts from Occurrences contexts statements.
# Dimensional contexts: Events (attributes). Order relations assertions
# by context occurrences hierarchy domain / range, set / superset
# attributes relations."
$old1 = "Dimensional statements. Occurrences contexts: Events. Dimensional contexts: Context occurrences hierarchy for order relations assertions."
$new1 = "Dimensional contexts: Contexts from Occurrences contexts statements. Dimensional contexts: Events (attributes). Order relations assertions by context occurrences hierarchy domain / range, set / superset attributes relations."

# --- Edit 2 -----------------------------------------------------------
# "(Mapping, Mapping super / parent / dimension, Kind unit, Role
# measure);" ->
# "Event (Dimensional context attributes): (Mapping / unit / class,
# Mapping super / parent / dimension / metaclass, Kind unit / measure /
# occurrence, Role measure / value instance);"
$old2 = "(Mapping, Mapping super / parent / dimension, Kind unit, Role measure);"
$new2 = "Event (Dimensional context attributes): (Mapping / unit / class, Mapping super / parent / dimension / metaclass, Kind unit / measure / occurrence, Role measure / value instance);"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $old1) {
        $p.Range.Text = $new1
    } elseif ($t -eq $old2) {
        $p.Range.Text = $new2
    }
}
